$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new "Given I am a role of a user" rows, one within each of
#     the last two scenario blocks, to make room for the missing
#     "Scenario:" header lines. ---
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(39).Insert()

# --- Fill in the newly inserted rows and fix the header text that shifted ---
$ws.Range("B34").Value = "Scenario: User sends photo"
$ws.Range("B35").Value = "Given I am a role of a user"

$ws.Range("B39").Value = "Scenario: User has questions"
$ws.Range("B40").Value = "Given I am a role of a user"

# --- Append a new trailing row 43 (just a styled, empty A cell) ---
$ws.Range("A43").Font.Bold = $true

# --- Bold column A across the whole CAP-44 / CAP-40 / CAP-39 / CAP-38 block ---
$ws.Range("A23:A43").Font.Bold = $true
